$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("status_label"), shifting NCTId/eudraCT/... etc. one column to the right.
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "status_label"

# Fill "status_label" = "rouge" for every data row (2-17), mirroring the "statut" = red square rows.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = "rouge"
}

# Rows 14-16 were also reordered (re-sorted) as part of this edit: the data that used to live in
# rows 15 and 16 now lives in rows 14 and 15, and the data that used to live in row 14 now lives in
# row 16 (a 3-way cyclic rotation). Columns C:G hold NCTId, eudraCT, completion_year,
# clinical_trial_title, acronym. We use Range.Copy (instead of re-typing the values) so that cell
# types - e.g. "2022" staying a text value rather than turning into a number - are preserved exactly
# like the other, untouched rows. A scratch row far below the data is used to stage row 14's values
# while the rotation happens, and is cleared again afterwards so it leaves no trace in the sheet.

$ws.Range("C14:G14").Copy($ws.Range("C100:G100"))
$ws.Range("C15:G15").Copy($ws.Range("C14:G14"))
$ws.Range("C16:G16").Copy($ws.Range("C15:G15"))
$ws.Range("C100:G100").Copy($ws.Range("C16:G16"))
$ws.Range("C100:G100").Clear()
